$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Unit fix in the raw measurement tables (ms -> s, and var/uncertainty scaled by 1/1e6) ---
# Row 2 (Cu2 250)
$ws.Range("D2").Value = 1.39484
$ws.Range("E2").Formula = "=0.001055/1000"
$ws.Range("F2").Value = 1.21551
$ws.Range("G2").Value = 0.0000002529

# Row 3 (Cu2 500)
$ws.Range("D3").Value = 1.0034
$ws.Range("E3").Formula = "=0.0004851/1000"
$ws.Range("F3").Value = 1.06644
$ws.Range("G3").Value = 0.0000002621

# Row 4 (Cu2 1000)
$ws.Range("D4").Value = 0.646849
$ws.Range("E4").Formula = "=71.54/1000"
$ws.Range("F4").Value = 0.748404
$ws.Range("G4").Value = 0.0000001937

# Row 5 (Cu2 2000)
$ws.Range("D5").Value = 0.431268
$ws.Range("E5").Formula = "=0.0002906/1000"
$ws.Range("F5").Value = 0.34183
$ws.Range("G5").Value = 0.0000001228

# Row 8 (Wasser)
$ws.Range("D8").Value = 2.19946
$ws.Range("E8").Value = 0.000003027
$ws.Range("F8").Value = 1.90106
$ws.Range("G8").Value = 0.08983

# Row 12 (Mn 2 25)
$ws.Range("D12").Value = 1.17828
$ws.Range("E12").Value = 0.0000009801
$ws.Range("F12").Value = 0.548337
$ws.Range("G12").Value = 0.0000001258

# Row 13 (Mn 2 50)
$ws.Range("D13").Value = 0.725857
$ws.Range("E13").Value = 0.0000006027
$ws.Range("F13").Value = 0.279858
$ws.Range("G13").Formula = "=8.179*10^(-8)"

# Row 14 (Mn 2 100)
$ws.Range("D14").Value = 0.316085
$ws.Range("E14").Value = 0.0000003079
$ws.Range("F14").Value = 0.170996
$ws.Range("G14").Value = 0.0000001182

# Row 15 (Mn 2 200)
$ws.Range("D15").Value = 0.180244
$ws.Range("E15").Value = 0.0000001274
$ws.Range("F15").Value = 0.0691512
$ws.Range("G15").Formula = "=5.47*10^(-8)"

# --- 2. Append a corrected r1/r2 summary block (rows 33:35), duplicating the
#        layout of the existing block in rows 29:31 with fixed-up numbers ---
$fmtSrc = $ws.Range("A29:M31")
$fmtSrc.Copy()
$ws.Range("A33").PasteSpecial(-4122)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M")
foreach ($c in $cols) {
    $ws.Range($c + "33").Value = $ws.Range($c + "29").Value2
}

$ws.Range("A34").Value = $ws.Range("A30").Value2
$ws.Range("B34").Value = 0.45395
$ws.Range("C34").Value = 0.03059
$ws.Range("D34").Value = 0.543404
$ws.Range("E34").Value = 0.07051
$ws.Range("F34").Formula = "=1/D34"
$ws.Range("G34").Formula = "=(1/D34^2)*E34"
$ws.Range("H34").Value = 0.616607
$ws.Range("I34").Value = 0.08373
$ws.Range("J34").Value = 0.349363
$ws.Range("K34").Value = 0.193
$ws.Range("L34").Formula = "=1/J34"
$ws.Range("M34").Formula = "=(1/J34^2)*K34"

$ws.Range("A35").Value = $ws.Range("A31").Value2
$ws.Range("B35").Value = 13.6523
$ws.Range("C35").Value = 0.8379
$ws.Range("D35").Value = 0.174718
$ws.Range("E35").Value = 0.1931
$ws.Range("F35").Formula = "=1/D35"
$ws.Range("G35").Formula = "=(1/D35^2)*E35"
$ws.Range("H35").Value = 35.9327
$ws.Range("I35").Value = 3.136
$ws.Range("J35").Value = -0.310853
$ws.Range("K35").Value = 0.7228
$ws.Range("L35").Formula = "=1/J35"
$ws.Range("M35").Formula = "=(1/J35^2)*K35"

# --- 3. Column E now holds numbers as wide as column D, so it picks up the same best-fit width ---
$ws.Columns("E").ColumnWidth = 11.17

# --- 4. Restore the view state (cursor parked on the new total row) ---
$ws.Range("F26").Select() | Out-Null
$ws.Range("M35").Select() | Out-Null
